$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update validation report values to reflect the latest pipeline run
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = -2

$ws.Range("E3").Value = 1

$ws.Range("D4").Value = -5

$ws.Range("B6").Value = 4

# Update the active selection on the sheet
$ws.Range("G8").Select()
